# Refresh the cryptos price table to the latest scrape values.
# Each entry is a cell reference -> new display text pulled from the
# source scrape. A few "Price" (column D) values are plain
# numeric-looking strings (e.g. "12.75"); a bare Value assignment of
# those would make Excel coerce the cell into a real number, which
# does not match the source (those columns hold formatted price
# text, some using "." as a thousands separator, e.g. "70.326.69").
# For such cells we force Text formatting before the write, then
# drop back to the default "Normal" style afterwards so no visible
# formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "70.326.69" }
    @{ Cell = "E2"; Value = "  +0.89%  " }
    @{ Cell = "D3"; Value = "3.565.90" }
    @{ Cell = "E3"; Value = "  +1.47%  " }
    @{ Cell = "E4"; Value = "  +0.03%  " }
    @{ Cell = "D5"; Value = "608.15" }
    @{ Cell = "E5"; Value = "  +3.55%  " }
    @{ Cell = "D6"; Value = "187.52" }
    @{ Cell = "E6"; Value = "  +2.04%  " }
    @{ Cell = "D7"; Value = "3.562.21" }
    @{ Cell = "E7"; Value = "  +1.68%  " }
    @{ Cell = "D8"; Value = "0.621" }
    @{ Cell = "E8"; Value = "  +1.46%  " }
    @{ Cell = "E9"; Value = "  -0.02%  " }
    @{ Cell = "D10"; Value = "0.213" }
    @{ Cell = "E10"; Value = "  +8.34%  " }
    @{ Cell = "D11"; Value = "0.647" }
    @{ Cell = "E11"; Value = "  +0.62%  " }
    @{ Cell = "D12"; Value = "53.98" }
    @{ Cell = "E12"; Value = "  -0.07%  " }
    @{ Cell = "D13"; Value = "0.0000309" }
    @{ Cell = "E13"; Value = "  +1.89%  " }
    @{ Cell = "D14"; Value = "9.41" }
    @{ Cell = "E14"; Value = "  -0.62%  " }
    @{ Cell = "D15"; Value = "4.136.24" }
    @{ Cell = "E15"; Value = "  +1.66%  " }
    @{ Cell = "D16"; Value = "70.430.32" }
    @{ Cell = "E16"; Value = "  +1.10%  " }
    @{ Cell = "B17"; Value = "WrappedEther" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" }
    @{ Cell = "D17"; Value = "3.573.91" }
    @{ Cell = "E17"; Value = "  +1.57%  " }
    @{ Cell = "B18"; Value = "Uniswap" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni" }
    @{ Cell = "D18"; Value = "12.75" }
    @{ Cell = "E18"; Value = "  +3.41%  " }
    @{ Cell = "D19"; Value = "18.98" }
    @{ Cell = "E19"; Value = "  -1.67%  " }
    @{ Cell = "D20"; Value = "573.66" }
    @{ Cell = "E20"; Value = "  +7.47%  " }
    @{ Cell = "D21"; Value = "0.120" }
    @{ Cell = "E21"; Value = "  +0.84%  " }
    @{ Cell = "D22"; Value = "0.994" }
    @{ Cell = "E22"; Value = "  -1.64%  " }
    @{ Cell = "D23"; Value = "17.48" }
    @{ Cell = "E23"; Value = "  -4.06%  " }
    @{ Cell = "D24"; Value = "4.78" }
    @{ Cell = "E24"; Value = "  +4.10%  " }
    @{ Cell = "D25"; Value = "4.95" }
    @{ Cell = "E25"; Value = "  +2.59%  " }
    @{ Cell = "D26"; Value = "94.22" }
    @{ Cell = "E26"; Value = "  -1.43%  " }
    @{ Cell = "D27"; Value = "2.93" }
    @{ Cell = "E27"; Value = "  -1.40%  " }
    @{ Cell = "D28"; Value = "10.94" }
    @{ Cell = "E28"; Value = "  -1.18%  " }
    @{ Cell = "D29"; Value = "9.39" }
    @{ Cell = "E29"; Value = "  +3.50%  " }
    @{ Cell = "D30"; Value = "32.28" }
    @{ Cell = "E30"; Value = "  +0.57%  " }
    @{ Cell = "D31"; Value = "7.07" }
    @{ Cell = "E31"; Value = "  -3.06%  " }
    @{ Cell = "D32"; Value = "12.21" }
    @{ Cell = "E32"; Value = "  -1.52%  " }
    @{ Cell = "B33"; Value = "OKB" }
    @{ Cell = "C33"; Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb" }
    @{ Cell = "D33"; Value = "64.49" }
    @{ Cell = "E33"; Value = "  +0.87%  " }
    @{ Cell = "B34"; Value = "Hedera" }
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar" }
    @{ Cell = "D34"; Value = "0.115" }
    @{ Cell = "E34"; Value = "  +1.46%  " }
    @{ Cell = "D35"; Value = "3.73" }
    @{ Cell = "E35"; Value = "  +20.45%  " }
    @{ Cell = "D36"; Value = "3.19" }
    @{ Cell = "E36"; Value = "  +2.69%  " }
    @{ Cell = "D37"; Value = "3.761.23" }
    @{ Cell = "E37"; Value = "  +12.17%  " }
    @{ Cell = "D38"; Value = "0.405" }
    @{ Cell = "E38"; Value = "  -0.57%  " }
    @{ Cell = "D39"; Value = "37.91" }
    @{ Cell = "E39"; Value = "  -0.65%  " }
    @{ Cell = "D40"; Value = "521.43" }
    @{ Cell = "E40"; Value = "  -4.54%  " }
    @{ Cell = "E41"; Value = "  +0.09%  " }
    @{ Cell = "D42"; Value = "0.0₃0782" }
    @{ Cell = "E42"; Value = "  +2.91%  " }
    @{ Cell = "D43"; Value = "3.54" }
    @{ Cell = "E43"; Value = "  +3.38%  " }
    @{ Cell = "D44"; Value = "0.138" }
    @{ Cell = "E44"; Value = "  +2.57%  " }
    @{ Cell = "D45"; Value = "0.0455" }
    @{ Cell = "E45"; Value = "  +3.84%  " }
    @{ Cell = "D46"; Value = "2.95" }
    @{ Cell = "E46"; Value = "  -0.55%  " }
    @{ Cell = "D47"; Value = "3.47" }
    @{ Cell = "E47"; Value = "  -0.98%  " }
    @{ Cell = "E48"; Value = "  +3.63%  " }
    @{ Cell = "D49"; Value = "9.19" }
    @{ Cell = "E49"; Value = "  +2.71%  " }
    @{ Cell = "E50"; Value = "  +0.36%  " }
    @{ Cell = "D51"; Value = "1.42" }
    @{ Cell = "E51"; Value = "  +6.66%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $looksNumeric = $u.Value -match "^-?\d+(\.\d+)?$"
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
